# Applies the "Add files via upload" revision to the Events/Employee workbook.
#
# Summary of the change (derived from the OOXML diff):
#   - On the "Events" sheet a new column is inserted after "Event" (B) and
#     before "Hall" (old C). It is titled "EventRanking" and holds a small
#     integer rank per event (1, 2, 3, 5, 7, 10).
#   - Everything that used to live in columns C..I (Hall, Skillset1,
#     Skillset2, Employees, Date, Shift begins, Shifts ends) shifts one
#     column to the right, landing in D..J.
#   - The last two headers are renamed:
#       "Shift begins" -> "ShiftBegins"
#       "Shifts ends"  -> "ShiftsEnds"
#   - The "Employee" sheet itself is untouched content-wise.
#   - Minor view cosmetics: the Events sheet is shown at 91% zoom with
#     C8 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- Insert the new "EventRanking" column between Event (B) and Hall (C) ---
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = 11

# --- Header row ---
$ws.Cells.Item(1, 3).Value = "EventRanking"

# Columns H..I (Shift begins / Shifts ends) moved to I..J after the insert;
# rename them to match the new naming convention.
$ws.Cells.Item(1, 9).Value = "ShiftBegins"
$ws.Cells.Item(1, 10).Value = "ShiftsEnds"

# --- Fill in the EventRanking values for the six events ---
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(6, 3).Value = 7
$ws.Cells.Item(7, 3).Value = 10

# --- View cosmetics to mirror the saved workbook state ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 91
$ws.Range("C8").Select()
